# SOSPaper / Tables_SOS.xlsx - "new fits and results"
#
# The "4) Model fits" sheet gets a second, more-precise results table
# (columns F:H, mirroring A:C) plus a small "updated" note in J2.
# A new shared string is introduced for that note.

$wb  = $excel.ActiveWorkbook
$ws4 = $wb.Worksheets.Item("4) Model fits")
$ws2 = $wb.Worksheets.Item("2) Model parameters")

# ---------------------------------------------------------------------
# 1) Bring over the cell formatting for the new F:J columns by copying
#    it from the existing, equivalently-styled cells on the same sheet
#    (and, for the one style that doesn't already exist on this sheet,
#    from "2) Model parameters"!A26).  We copy formats only, then set
#    the values/text afterwards so we don't disturb anything else.
# ---------------------------------------------------------------------

# Row 1 - table title, mirrored into F1
$ws4.Range("A1").Copy()
$ws4.Range("F1").PasteSpecial(-4122)

# Row 2 - header row, mirrored into F2:H2 ; J2 gets the "updated" note
$ws4.Range("A2").Copy()
$ws4.Range("F2").PasteSpecial(-4122)
$ws4.Range("B2").Copy()
$ws4.Range("G2").PasteSpecial(-4122)
$ws4.Range("C2").Copy()
$ws4.Range("H2").PasteSpecial(-4122)
$ws2.Range("A26").Copy()
$ws4.Range("J2").PasteSpecial(-4122)

# Rows 3-6 - lake name label mirrored into column F (plain style like A3:A6);
# the new RMSE/NSE numbers in G/H keep the workbook's default (unstyled) look
$ws4.Range("A3").Copy()
$ws4.Range("F3").PasteSpecial(-4122)
$ws4.Range("A4").Copy()
$ws4.Range("F4").PasteSpecial(-4122)
$ws4.Range("A5").Copy()
$ws4.Range("F5").PasteSpecial(-4122)
$ws4.Range("A6").Copy()
$ws4.Range("F6").PasteSpecial(-4122)

# Row 7 - bottom border / thick-bot row style
$ws4.Range("A7").Copy()
$ws4.Range("F7").PasteSpecial(-4122)

# Rows 8-9 - footnote label style
$ws4.Range("A8").Copy()
$ws4.Range("F8").PasteSpecial(-4122)
$ws4.Range("A9").Copy()
$ws4.Range("F9").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2) Fill in the text/values for the mirrored table + new note
# ---------------------------------------------------------------------

$ws4.Range("F1").Value = "Table 4. Model goodness of fit"

$ws4.Range("F2").Value = "Lake"
$ws4.Range("G2").Value = "RMSE*"
$ws4.Range("H2").Value = "NSE**"
$ws4.Range("J2").Value = "<-- Updated May 8 2018 HD"

$ws4.Range("F3").Value = "Harp"
$ws4.Range("G3").Value = 1.2488489306156501
$ws4.Range("H3").Value = 0.85389351218815102

$ws4.Range("F4").Value = "Monona"
$ws4.Range("G4").Value = 1.4641714180909799
$ws4.Range("H4").Value = 0.69837503040359505

$ws4.Range("F5").Value = "Toolik"
$ws4.Range("G5").Value = 1.2123881443494
$ws4.Range("H5").Value = 0.79254393893729602

$ws4.Range("F6").Value = "Trout"
$ws4.Range("G6").Value = 0.80209904665356302
$ws4.Range("H6").Value = 0.95158047733254303

$ws4.Range("F7").Value = "Vanern"
$ws4.Range("G7").Value = 0.72408838354316896
$ws4.Range("H7").Value = 0.965139555519964

$ws4.Range("F8").Value = "* root mean square error (mg/L)"
$ws4.Range("F9").Value = "** Nash-Sutcliffe efficiency"

# ---------------------------------------------------------------------
# 3) Window / selection state: "4) Model fits" becomes the active sheet
#    (it picks up tabSelected), with E3 selected.
# ---------------------------------------------------------------------

$ws4.Activate()
$ws4.Range("E3").Select()
